$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Capture the existing values of the 4 columns that will shift right (K..N -> M..P) ---
$K1 = $ws.Range("K1").Value2; $K2 = $ws.Range("K2").Value2; $K3 = $ws.Range("K3").Value2
$L1 = $ws.Range("L1").Value2; $L2 = $ws.Range("L2").Value2; $L3 = $ws.Range("L3").Value2
$M1 = $ws.Range("M1").Value2; $M2 = $ws.Range("M2").Value2; $M3 = $ws.Range("M3").Value2
$N1 = $ws.Range("N1").Value2; $N2 = $ws.Range("N2").Value2; $N3 = $ws.Range("N3").Value2

# --- Grow the table to the right by 2 columns (B2:N3 -> B2:P3) ---
$lo.Resize($ws.Range("B2:P3"))

# --- Re-home the old ItemSlotMax/BallStackMax/GoldGain/DefaultBallDataId data two columns right ---
$ws.Range("M1").Value = $K1; $ws.Range("M2").Value = $K2; $ws.Range("M3").Value = $K3
$ws.Range("N1").Value = $L1; $ws.Range("N2").Value = $L2; $ws.Range("N3").Value = $L3
$ws.Range("O1").Value = $M1; $ws.Range("O2").Value = $M2; $ws.Range("O3").Value = $M3
$ws.Range("P1").Value = $N1; $ws.Range("P2").Value = $N2; $ws.Range("P3").Value = $N3

# --- Populate the two new columns: DodgeDurationTime (K) and DodgeAddForce (L) ---
$ws.Range("L1").Value = "float"
$ws.Range("K1").Value = "float"
$ws.Range("L2").Value = "DodgeAddForce"
$ws.Range("K2").Value = "DodgeDurationTime"
$ws.Range("K3").Value = 0.3
$ws.Range("L3").Value = 12

# --- Fix up the ListColumn names so the table definition matches the new header text ---
$lo.ListColumns.Item(10).Name = "DodgeDurationTime"
$lo.ListColumns.Item(11).Name = "DodgeAddForce"
$lo.ListColumns.Item(12).Name = "ItemSlotMax"
$lo.ListColumns.Item(13).Name = "BallStackMax"
$lo.ListColumns.Item(14).Name = "GoldGain"
$lo.ListColumns.Item(15).Name = "DefaultBallDataId"

# --- Player MovementSpd (F3) tuning: 300 -> 7 ---
$ws.Range("F3").Value = 7

# --- Column widths: K:L share the MovementSpd-style width, M:P keep their old bestFit widths, Q gets a small pad column ---
$ws.Columns("K:L").ColumnWidth = 22.75
$ws.Columns("Q:Q").ColumnWidth = 9.14

# --- Restore selection roughly where the author left off ---
$ws.Range("L8").Select()
